# Add support for the "call_and_response" long-tone singing paradigm
# alongside the existing "sing_along" one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New dictionary rows (146-150) -----------------------------------------
# Column A: new translation keys (these become new shared strings).
$ws.Range("A146").Value = "long_tone_title_call_and_response"
$ws.Range("A147").Value = "long_tone_text_call_and_response"
$ws.Range("A148").Value = "long_tone_instruction_call_and_response"
$ws.Range("A149").Value = "long_tone_instruction_call_and_response_2"
$ws.Range("A150").Value = "long_tone_instruction_call_and_response_3"

# Columns B/C/D: the (English) text values repeated across the three
# language columns, matching the pattern used by the existing rows.
$ws.Range("B146").Value = "Sing Back The Note After Your Hear It"
$ws.Range("C146").Value = "Sing Back The Note After Your Hear It"
$ws.Range("D146").Value = "Sing Back The Note After Your Hear It"

$ws.Range("B147").Value = "Sing back the note for 5 seconds after you hear it."
$ws.Range("C147").Value = "Sing back the note for 5 seconds after you hear it."
$ws.Range("D147").Value = "Sing back the note for 5 seconds after you hear it."

$ws.Range("B148").Value = "When you click the Play button in the next set of trials, you will hear a 5-second note."
$ws.Range("C148").Value = "When you click the Play button in the next set of trials, you will hear a 5-second note."
$ws.Range("D148").Value = "When you click the Play button in the next set of trials, you will hear a 5-second note."

$ws.Range("B149").Value = "Please try and sing the exact same note and hold after you hear it."
$ws.Range("C149").Value = "Please try and sing the exact same note and hold after you hear it."
$ws.Range("D149").Value = "Please try and sing the exact same note and hold after you hear it."

$ws.Range("B150").Value = "If the note seems out of your voice range, sing a note that is in your range that best matches the note, for example, in a different octave. "
$ws.Range("C150").Value = "If the note seems out of your voice range, sing a note that is in your range that best matches the note, for example, in a different octave. "
$ws.Range("D150").Value = "If the note seems out of your voice range, sing a note that is in your range that best matches the note, for example, in a different octave. "

# --- Formatting: reuse the styles from the equivalent "sing_along" rows ----
# Column A (rows 146-150) mirrors A83's style (key column formatting).
$ws.Range("A83").Copy()
$ws.Range("A146:A150").PasteSpecial(-4122)

# Columns B:D for rows 146-149 mirror B83:D83's style.
$ws.Range("B83:D83").Copy()
$ws.Range("B146:D149").PasteSpecial(-4122)

# Columns B:D for row 150 mirror B7's style (left/vcenter aligned variant).
$ws.Range("B7").Copy()
$ws.Range("B150:D150").PasteSpecial(-4122)

# --- Sheet-level formatting updates -----------------------------------------
# Widen column A to fit the new, longer keys.
$ws.Range("A1:A150").ColumnWidth = 56.142857142857146

# --- View state updates (best effort) ---------------------------------------
$null = $ws.Range("D152").Select()
